# Apply the statement_14.xlsx update:
# - new card holder name / card number
# - new transaction period (22.07.2025 - 05.08.2025) replacing old one (01.04.2024 - 18.04.2024)
# - one fewer transaction row (row 10 becomes blank, matching the already-blank row 11 style)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: card holder name / card number
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long digit string that must stay text (as in the original file),
# not get auto-converted to a number: force text format before assigning,
# then re-apply the original cell formatting (copied from a same-styled
# neighbour) so the cell keeps its original style index.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("D3").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 22.07.2025"

# Transaction row 6
$ws.Range("B6").Value = "23.07."
$ws.Range("C6").Value = "24.07."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 58924974"
$ws.Range("E6").Value = "85,09-"

# Transaction row 7
$ws.Range("B7").Value = "26.07."
$ws.Range("C7").Value = "27.07."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-145769"
$ws.Range("E7").Value = "54,20-"

# Transaction row 8
$ws.Range("B8").Value = "28.07."
$ws.Range("C8").Value = "29.07."
$ws.Range("D8").Value = "KARTENZ./28.07 LIDL RO"
$ws.Range("E8").Value = "39,70-"

# Transaction row 9
$ws.Range("B9").Value = "01.08."
$ws.Range("C9").Value = "02.08."
$ws.Range("D9").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E9").Value = "25,30-"

# Row 10 no longer holds a transaction - clear its contents, and copy the
# formatting already used by the (blank) row 11 onto E10's cell so the
# "amount" column keeps the correct right/vertical-centred/wrap styling
# used for blank rows in this table.
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("E11").Copy()
$ws.Range("E10").PasteSpecial(-4122)

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 05.08.2025"
$ws.Range("E12").Value = "204,29-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 14.08.2025"
